# Restored from revision of admin on 06/24/2020 07:23:36 AM.TEST Author: admin. Type: SAVE.
# The only functional change in this revision is the value of cell C10 on the
# "Rules" sheet, which goes from 18 to 1 (the R30 rule's "Integer min" bound).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
